# Update the "Minimum Renewable Percentage" scenario table.
# The 0.65 "renewable floor" block that used to apply to the historic
# years 2010-2017 (columns E:L) now applies to the future years
# 2024-2050 (columns S:AS) instead; 2010-2017 drops to 0. The years
# 2018-2023 (columns M:R) stay at 0, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-18: columns E:L (2010-2017) -> 0
$ws.Range("E2:L18").Value = 0

# Data rows 2-18: columns S:AS (2024-2050) -> 0.65
$ws.Range("S2:AS18").Value = 0.65

# Match the author's final selection in the sheet.
$null = $ws.Range("S2:AS18").Select()
